$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Paragraph 1: "Sprint Reviews and Retrospectives:" -> "Sprint Retrospectives:"
$para1 = $tr.Paragraphs(1)
$para1.Runs(1).Text = "Sprint Retrospectives:"

# Paragraph 2: update first run and the trailing run, keep the bold "reward to the team" run intact
$para2 = $tr.Paragraphs(2)
$para2.Runs(1).Text = "The Sprint Retrospective is an informal meeting that should feel as a "
$para2.Runs(3).Text = " for the completed Sprint. This reward is a privilege that allows them to retrospectively learn from experiences."

# Paragraph 3: rewording
$para3 = $tr.Paragraphs(3)
$para3.Runs(1).Text = "Unless you face the extreme case of the meeting going over two hours, you do not want to limit its duration. Therefore, do not have it standing. The team must feel comfortable to spend whichever reasonable amount of time they want discussing how they did things."
